# Natmi following Dr Hou advice
#
# Re-run of the NATMI ligand-receptor scoring for Sema5a-Met with the
# "sCs" / "ECs" / "FAPs" sending-cluster set (per Dr Hou's advice). The
# target-cluster set per sending cluster also grows from 3 to 4
# (ECs, FAPs, M2, sCs), so the table grows from 8 data rows (A2:T9) to
# 12 data rows (A2:T13) and every numeric column is recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema5a"
$ws.Cells.Item(2,3).Value = "Met"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.5290683333333334
$ws.Cells.Item(2,8).Value = 1.587205
$ws.Cells.Item(2,9).Value = 0.01267142171338989
$ws.Cells.Item(2,10).Value = 0.01267142171338989
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.655851666666667
$ws.Cells.Item(2,14).Value = 4.967555
$ws.Cells.Item(2,15).Value = 0.03628213169899143
$ws.Cells.Item(2,16).Value = 0.03628213169899143
$ws.Cells.Item(2,17).Value = 0.8760586815305557
$ws.Cells.Item(2,18).Value = 7.884528133775
$ws.Cells.Item(2,19).Value = 0.0004597461914186715
$ws.Cells.Item(2,20).Value = 0.0004597461914186715
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema5a"
$ws.Cells.Item(3,3).Value = "Met"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.5290683333333334
$ws.Cells.Item(3,8).Value = 1.587205
$ws.Cells.Item(3,9).Value = 0.01267142171338989
$ws.Cells.Item(3,10).Value = 0.01267142171338989
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.416382666666667
$ws.Cells.Item(3,14).Value = 4.249148
$ws.Cells.Item(3,15).Value = 0.03103501568568562
$ws.Cells.Item(3,16).Value = 0.03103501568568562
$ws.Cells.Item(3,17).Value = 0.7493632168155555
$ws.Cells.Item(3,18).Value = 6.74426895134
$ws.Cells.Item(3,19).Value = 0.0003932577716349924
$ws.Cells.Item(3,20).Value = 0.0003932577716349925
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema5a"
$ws.Cells.Item(4,3).Value = "Met"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.5290683333333334
$ws.Cells.Item(4,8).Value = 1.587205
$ws.Cells.Item(4,9).Value = 0.01267142171338989
$ws.Cells.Item(4,10).Value = 0.01267142171338989
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.385314999999999
$ws.Cells.Item(4,14).Value = 16.155945
$ws.Cells.Item(4,15).Value = 0.1180001276707882
$ws.Cells.Item(4,16).Value = 0.1180001276707882
$ws.Cells.Item(4,17).Value = 2.849199631525
$ws.Cells.Item(4,18).Value = 25.642796683725
$ws.Cells.Item(4,19).Value = 0.001495229379950404
$ws.Cells.Item(4,20).Value = 0.001495229379950404
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Sema5a"
$ws.Cells.Item(5,3).Value = "Met"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.5290683333333334
$ws.Cells.Item(5,8).Value = 1.587205
$ws.Cells.Item(5,9).Value = 0.01267142171338989
$ws.Cells.Item(5,10).Value = 0.01267142171338989
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 37.180664
$ws.Cells.Item(5,14).Value = 111.541992
$ws.Cells.Item(5,15).Value = 0.8146827249445348
$ws.Cells.Item(5,16).Value = 0.8146827249445348
$ws.Cells.Item(5,17).Value = 19.67111193470667
$ws.Cells.Item(5,18).Value = 177.04000741236
$ws.Cells.Item(5,19).Value = 0.01032318837038582
$ws.Cells.Item(5,20).Value = 0.01032318837038582
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sema5a"
$ws.Cells.Item(6,3).Value = "Met"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 32.04971933333334
$ws.Cells.Item(6,8).Value = 96.149158
$ws.Cells.Item(6,9).Value = 0.7676050216609417
$ws.Cells.Item(6,10).Value = 0.7676050216609416
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.655851666666667
$ws.Cells.Item(6,14).Value = 4.967555
$ws.Cells.Item(6,15).Value = 0.03628213169899143
$ws.Cells.Item(6,16).Value = 0.03628213169899143
$ws.Cells.Item(6,17).Value = 53.0695811742989
$ws.Cells.Item(6,18).Value = 477.62623056869
$ws.Cells.Item(6,19).Value = 0.02785034648870946
$ws.Cells.Item(6,20).Value = 0.02785034648870945
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema5a"
$ws.Cells.Item(7,3).Value = "Met"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 32.04971933333334
$ws.Cells.Item(7,8).Value = 96.149158
$ws.Cells.Item(7,9).Value = 0.7676050216609417
$ws.Cells.Item(7,10).Value = 0.7676050216609416
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.416382666666667
$ws.Cells.Item(7,14).Value = 4.249148
$ws.Cells.Item(7,15).Value = 0.03103501568568562
$ws.Cells.Item(7,16).Value = 0.03103501568568562
$ws.Cells.Item(7,17).Value = 45.39466693526489
$ws.Cells.Item(7,18).Value = 408.552002417384
$ws.Cells.Item(7,19).Value = 0.02382263388765837
$ws.Cells.Item(7,20).Value = 0.02382263388765837
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Sema5a"
$ws.Cells.Item(8,3).Value = "Met"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 32.04971933333334
$ws.Cells.Item(8,8).Value = 96.149158
$ws.Cells.Item(8,9).Value = 0.7676050216609417
$ws.Cells.Item(8,10).Value = 0.7676050216609416
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.385314999999999
$ws.Cells.Item(8,14).Value = 16.155945
$ws.Cells.Item(8,15).Value = 0.1180001276707882
$ws.Cells.Item(8,16).Value = 0.1180001276707882
$ws.Cells.Item(8,17).Value = 172.59783427159
$ws.Cells.Item(8,18).Value = 1553.38050844431
$ws.Cells.Item(8,19).Value = 0.09057749055672923
$ws.Cells.Item(8,20).Value = 0.09057749055672923
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Sema5a"
$ws.Cells.Item(9,3).Value = "Met"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 32.04971933333334
$ws.Cells.Item(9,8).Value = 96.149158
$ws.Cells.Item(9,9).Value = 0.7676050216609417
$ws.Cells.Item(9,10).Value = 0.7676050216609416
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 37.180664
$ws.Cells.Item(9,14).Value = 111.541992
$ws.Cells.Item(9,15).Value = 0.8146827249445348
$ws.Cells.Item(9,16).Value = 0.8146827249445348
$ws.Cells.Item(9,17).Value = 1191.629845826971
$ws.Cells.Item(9,18).Value = 10724.66861244274
$ws.Cells.Item(9,19).Value = 0.6253545507278446
$ws.Cells.Item(9,20).Value = 0.6253545507278445
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Sema5a"
$ws.Cells.Item(10,3).Value = "Met"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 9.174090999999999
$ws.Cells.Item(10,8).Value = 27.522273
$ws.Cells.Item(10,9).Value = 0.2197235566256685
$ws.Cells.Item(10,10).Value = 0.2197235566256685
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.655851666666667
$ws.Cells.Item(10,14).Value = 4.967555
$ws.Cells.Item(10,15).Value = 0.03628213169899143
$ws.Cells.Item(10,16).Value = 0.03628213169899143
$ws.Cells.Item(10,17).Value = 15.19093387250166
$ws.Cells.Item(10,18).Value = 136.718404852515
$ws.Cells.Item(10,19).Value = 0.007972039018863305
$ws.Cells.Item(10,20).Value = 0.007972039018863305
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Sema5a"
$ws.Cells.Item(11,3).Value = "Met"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 9.174090999999999
$ws.Cells.Item(11,8).Value = 27.522273
$ws.Cells.Item(11,9).Value = 0.2197235566256685
$ws.Cells.Item(11,10).Value = 0.2197235566256685
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.416382666666667
$ws.Cells.Item(11,14).Value = 4.249148
$ws.Cells.Item(11,15).Value = 0.03103501568568562
$ws.Cells.Item(11,16).Value = 0.03103501568568562
$ws.Cells.Item(11,17).Value = 12.99402347482267
$ws.Cells.Item(11,18).Value = 116.946211273404
$ws.Cells.Item(11,19).Value = 0.006819124026392253
$ws.Cells.Item(11,20).Value = 0.006819124026392255
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Sema5a"
$ws.Cells.Item(12,3).Value = "Met"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 9.174090999999999
$ws.Cells.Item(12,8).Value = 27.522273
$ws.Cells.Item(12,9).Value = 0.2197235566256685
$ws.Cells.Item(12,10).Value = 0.2197235566256685
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 5.385314999999999
$ws.Cells.Item(12,14).Value = 16.155945
$ws.Cells.Item(12,15).Value = 0.1180001276707882
$ws.Cells.Item(12,16).Value = 0.1180001276707882
$ws.Cells.Item(12,17).Value = 49.40536987366499
$ws.Cells.Item(12,18).Value = 444.648328862985
$ws.Cells.Item(12,19).Value = 0.02592740773410853
$ws.Cells.Item(12,20).Value = 0.02592740773410854
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Sema5a"
$ws.Cells.Item(13,3).Value = "Met"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 9.174090999999999
$ws.Cells.Item(13,8).Value = 27.522273
$ws.Cells.Item(13,9).Value = 0.2197235566256685
$ws.Cells.Item(13,10).Value = 0.2197235566256685
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 37.180664
$ws.Cells.Item(13,14).Value = 111.541992
$ws.Cells.Item(13,15).Value = 0.8146827249445348
$ws.Cells.Item(13,16).Value = 0.8146827249445348
$ws.Cells.Item(13,17).Value = 341.098794976424
$ws.Cells.Item(13,18).Value = 3069.889154787816
$ws.Cells.Item(13,19).Value = 0.1790049858463044
$ws.Cells.Item(13,20).Value = 0.1790049858463044
